$d = $word.ActiveDocument

$replacements = @(
    @("159÷6=26, 3", "803÷5=160, 3"),
    @("341÷9=37, 8", "955÷4=238, 3"),
    @("414÷8=51, 6", "715÷7=102, 1"),
    @("855÷4=213, 3", "703÷5=140, 3"),
    @("730÷4=182, 2", "762÷4=190, 2"),
    @("423÷2=211, 1", "739÷8=92, 3"),
    @("110÷4=27, 2", "662÷6=110, 2"),
    @("862÷5=172, 2", "504÷9=56, 0"),
    @("109÷6=18, 1", "734÷5=146, 4"),
    @("991÷9=110, 1", "303÷9=33, 6"),
    @("643÷4=160, 3", "287÷5=57, 2"),
    @("165÷2=82, 1", "839÷4=209, 3"),
    @("722÷2=361, 0", "961÷4=240, 1"),
    @("281÷2=140, 1", "949÷3=316, 1"),
    @("816÷9=90, 6", "975÷5=195, 0"),
    @("805÷6=134, 1", "604÷3=201, 1"),
    @("862÷6=143, 4", "706÷2=353, 0"),
    @("963÷5=192, 3", "560÷7=80, 0"),
    @("534÷5=106, 4", "778÷9=86, 4"),
    @("244÷6=40, 4", "205÷8=25, 5"),
    @("308÷5=61, 3", "417÷8=52, 1"),
    @("887÷4=221, 3", "465÷9=51, 6"),
    @("864÷9=96, 0", "287÷8=35, 7"),
    @("786÷3=262, 0", "476÷5=95, 1"),
    @("712÷6=118, 4", "820÷5=164, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
